# This script applies a row-level reshuffle to the data rows (2-14) of the
# single worksheet in the workbook. Every data row (columns A:T) is an
# independent record; the edit re-orders these records across rows 2-14
# without changing the header row (row 1) or any column structure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row number -> source row number (both in the
# "before" snapshot). We snapshot every source row's full A:T values
# first, then write them back out in the new order, so that the
# read-before-write ordering never clobbers data we still need.
$rowMap = @{
    2  = 9
    3  = 7
    4  = 12
    5  = 4
    6  = 13
    7  = 14
    8  = 5
    9  = 8
    10 = 3
    11 = 10
    12 = 2
    13 = 11
    14 = 6
}

$firstCol = 1   # A
$lastCol  = 20  # T

# Snapshot all source rows (2-14) fully before writing anything back.
# Note: use Value2 (not Value) for multi-cell ranges, since Value does not
# reliably marshal a real 2-D array in this COM-interop runtime.
$snapshot = @{}
for ($r = 2; $r -le 14; $r++) {
    $rng = $ws.Range($ws.Cells.Item($r, $firstCol), $ws.Cells.Item($r, $lastCol))
    $snapshot[$r] = $rng.Value2
}

# Write the snapshot values to their new destination rows.
foreach ($destRow in ($rowMap.Keys | Sort-Object)) {
    $srcRow = $rowMap[$destRow]
    $values = $snapshot[$srcRow]
    $destRng = $ws.Range($ws.Cells.Item($destRow, $firstCol), $ws.Cells.Item($destRow, $lastCol))
    $destRng.Value2 = $values
}
